$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "2022-Q3" sheet before the existing "2022-Q1" sheet
#    (position 2), by duplicating "2022-Q1" so the header row / styles /
#    column widths all come along for free, then overwriting its data
#    row with the 2022-Q3 figures.
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item(2)
$q1Sheet.Copy($q1Sheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Keep these as text (matching how the source data is stored) by using
# a leading apostrophe so Excel doesn't coerce the numeric-looking
# strings into real numbers.
$q3Sheet.Cells.Item(2, 4).Value = "'4.34"
$q3Sheet.Cells.Item(2, 5).Value = "'94.11"
$q3Sheet.Cells.Item(2, 6).Value = "'2.49"
$q3Sheet.Cells.Item(2, 7).Value = "'0.1081"
$q3Sheet.Cells.Item(2, 8).Value = 10

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: a new row for 2022-Q3 is added at
#    the top of the data, and every existing row shifts down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

# Row 8 is brand new - clone the formatting of row 7's A cell (bold,
# bordered, centered) before filling in the 2020-Q4 values that used to
# live in row 7.
$totalSheet.Cells.Item(7, 1).Copy()
$totalSheet.Cells.Item(8, 1).PasteSpecial(-4122)

$totalSheet.Cells.Item(8, 1).Value = 6
$totalSheet.Cells.Item(8, 2).Value = "2020-Q4"
$totalSheet.Cells.Item(8, 3).Value = 25
$totalSheet.Cells.Item(8, 4).Value = 0.65

$totalSheet.Cells.Item(7, 1).Value = 5
$totalSheet.Cells.Item(7, 2).Value = "2021-Q1"
$totalSheet.Cells.Item(7, 3).Value = 5
$totalSheet.Cells.Item(7, 4).Value = 0.27

$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(6, 2).Value = "2021-Q2"
$totalSheet.Cells.Item(6, 3).Value = 5
$totalSheet.Cells.Item(6, 4).Value = 0.27

$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(5, 3).Value = 2
$totalSheet.Cells.Item(5, 4).Value = 0.28

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(4, 3).Value = 1
$totalSheet.Cells.Item(4, 4).Value = 0.21

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(3, 3).Value = 1
$totalSheet.Cells.Item(3, 4).Value = 0.17

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 1
$totalSheet.Cells.Item(2, 4).Value = 0.11
